$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Insert a brand-new order row at row 17 on sheet "订单" (pushes the three
#    existing rows 17-19 down to 18-20, dimension becomes A1:J20).
# ---------------------------------------------------------------------------
$ws1.Rows("17:17").Insert()

# ---------------------------------------------------------------------------
# 2. Fill in the new row 17 with the new "曹卓肺癌和癌旁组织对比分析" order.
# ---------------------------------------------------------------------------
$ws1.Range("A17").Value = 45258
$ws1.Range("B17").Value = 16
$ws1.Range("C17").Value = ""
$ws1.Range("D17").Value = "曹卓肺癌和癌旁组织对比分析"
$ws1.Range("E17").Value = ""
$ws1.Range("F17").Value = "黄礼闯"
$ws1.Range("G17").Value = 45244
$ws1.Range("H17").Value = "完成"
$ws1.Range("I17").Value = "肺癌和癌旁组织单细胞数据对比分析"
$ws1.Range("J17").Value = ""

# ---------------------------------------------------------------------------
# 3. The row insert shifted the old rows 17-19 down to 18-20 but left their
#    "序号" (B column) values untouched; renumber them sequentially.
# ---------------------------------------------------------------------------
$ws1.Range("B18").Value = 17
$ws1.Range("B19").Value = 18
$ws1.Range("B20").Value = 19

# ---------------------------------------------------------------------------
# 4. The "曹卓补充订单" row (now row 20) moved from "待完成" to "完成".
# ---------------------------------------------------------------------------
$ws1.Range("H20").Value = "完成"

# ---------------------------------------------------------------------------
# 5. Re-stamp every date cell (派发日期/日期 columns) with the date format so
#    they all share one refreshed style, matching the re-touched look of the
#    sheet after the new row was spliced in.
# ---------------------------------------------------------------------------
$ws1.Range("N1").Value = 1
$ws1.Range("N1").NumberFormat = "mm-dd-yy"
$dateFormat = $ws1.Range("N1").NumberFormat
$ws1.Range("N1").Clear()

$ws1.Range("A2:A20,G2:G20").NumberFormat = $dateFormat
$ws2.Range("A2:A4,G2:G4").NumberFormat = $dateFormat

Write-Host "done"
